$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'69.687.07"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.47%  "
$ws.Range("D3").Value = "'3.892.53"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.25%  "
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").Value = "'604.58"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.58%  "
$ws.Range("D6").Value = "'171.46"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.61%  "
$ws.Range("D7").Value = "'3.891.72"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.35%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("E9").Value = "  +1.18%  "
$ws.Range("E10").Value = "  +1.32%  "
$ws.Range("D11").Value = "'6.40"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.73%  "
$ws.Range("E12").Value = "  +1.69%  "
$ws.Range("D13").Value = "'0.0000256"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.85%  "
$ws.Range("D14").Value = "'38.26"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.50%  "
$ws.Range("D15").Value = "'4.549.78"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Value = "'3.899.87"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.96%  "
$ws.Range("D17").Value = "'69.674.96"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.22%  "
$ws.Range("D18").Value = "'18.72"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +9.16%  "
$ws.Range("D19").Value = "'7.61"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.79%  "
$ws.Range("E20").Value = "  -0.73%  "
$ws.Range("D21").Value = "'11.05"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.60%  "
$ws.Range("D22").Value = "'490.21"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.83%  "
$ws.Range("D23").Value = "'0.746"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.85%  "
$ws.Range("D24").Value = "'0.0000166"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.19%  "
$ws.Range("D25").Value = "'85.37"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.53%  "
$ws.Range("D26").Value = "'2.30"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.66%  "
$ws.Range("D27").Value = "'12.35"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.89%  "
$ws.Range("D28").Value = "'10.14"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.40%  "
$ws.Range("E29").Value = "  +0.20%  "
$ws.Range("E30").Value = "  +1.14%  "
$ws.Range("D31").Value = "'4.044.92"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.22%  "
$ws.Range("D32").Value = "'2.40"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.31%  "
$ws.Range("D33").Value = "'7.84"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.31%  "
$ws.Range("D34").Value = "'31.90"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.24%  "
$ws.Range("D35").Value = "'3.863.35"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.92%  "
$ws.Range("E36").Value = "  -0.23%  "
$ws.Range("D37").Value = "'3.43"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +15.13%  "
$ws.Range("D38").Value = "'6.12"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.91%  "
$ws.Range("E39").Value = "  +2.12%  "
$ws.Range("E40").Value = "  +0.73%  "
$ws.Range("D41").Value = "'0.999"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.06%  "
$ws.Range("D42").Value = "'0.327"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.69%  "
$ws.Range("D43").Value = "'2.09"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +5.13%  "
$ws.Range("D44").Value = "'435.28"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.64%  "
$ws.Range("E45").Value = "  -0.92%  "
$ws.Range("D46").Value = "'8.69"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.37%  "
$ws.Range("D48").Value = "'0.000275"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +20.86%  "
$ws.Range("D49").Value = "'0.0366"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.50%  "
$ws.Range("D50").Value = "'40.27"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.76%  "
$ws.Range("D51").Value = "'142.39"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.24%  "
